$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 1.094976487407548

# Row 3
$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 5.553084769722144

# Row 4
$ws.Range("B4").Value = 1.459612070389937
$ws.Range("C4").Value = 0.3127903958511391
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 3.074671312995807
